$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158. This shifts the existing rows 158-214
# down to 159-215, carrying their formatting (including the date style on
# column D) along with them.
$ws.Rows.Item(158).Insert()

# The newly inserted row 158 now duplicates the row that used to be at 158
# (now row 159) in terms of formatting. Populate it with the new record's
# values as described by the commit.
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44900
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112012
$ws.Range("G158").Value = "Espinaca"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 55
$ws.Range("K158").Value = 10000
$ws.Range("L158").Value = 10000
$ws.Range("M158").Value = 10000
$ws.Range("N158").Value = "$/docena de atados"
$ws.Range("O158").Value = "Región de La Araucanía"
$ws.Range("P158").Value = 3333
$ws.Range("Q158").Value = 3
$ws.Range("R158").Value = "Hortaliza"
